$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("D2").Value = "29.202.01"
$ws.Range("D3").Value = "1.828.17"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "'234.11"
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("D6").Value = "'0.5989"
$ws.Range("E6").Value = "  -4.45%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "'0.06970"
$ws.Range("E8").Value = "  -5.87%  "
$ws.Range("D9").Value = "'0.2754"
$ws.Range("E9").Value = "  -4.85%  "
$ws.Range("D10").Value = "'23.23"
$ws.Range("E10").Value = "  -6.75%  "
$ws.Range("D11").Value = "'0.07616"
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("D12").Value = "1.838.38"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").Value = "'4.752"
$ws.Range("E13").Value = "  -4.33%  "
$ws.Range("D14").Value = "'0.6260"
$ws.Range("E14").Value = "  -6.69%  "
$ws.Range("D15").Value = "'0.000009652"
$ws.Range("E15").Value = "  -6.81%  "
$ws.Range("D16").Value = "'78.27"
$ws.Range("E16").Value = "  -4.30%  "
$ws.Range("D17").Value = "28.847.42"
$ws.Range("E17").Value = "  -1.84%  "
$ws.Range("D18").Value = "'5.704"
$ws.Range("E18").Value = "  -8.87%  "
$ws.Range("D19").Value = "'220.68"
$ws.Range("E19").Value = "  -5.74%  "
$ws.Range("D20").Value = "'1.004"
$ws.Range("E21").Value = "  -6.32%  "
$ws.Range("D22").Value = "'6.844"
$ws.Range("E22").Value = "  -6.33%  "
$ws.Range("D23").Value = "'1.005"
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("D24").Value = "'155.42"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").Value = "'7.955"
$ws.Range("E25").Value = "  -6.18%  "
$ws.Range("D26").Value = "'0.1289"
$ws.Range("E26").Value = "  -4.21%  "
$ws.Range("D27").Value = "'16.52"
$ws.Range("E27").Value = "  -4.63%  "
$ws.Range("D28").Value = "'0.06514"
$ws.Range("E28").Value = "  -10.05%  "
$ws.Range("D29").Value = "'1.448"
$ws.Range("E29").Value = "  -2.69%  "
$ws.Range("D31").Value = "'3.831"
$ws.Range("E31").Value = "  -5.06%  "
$ws.Range("D32").Value = "'3.758"
$ws.Range("E32").Value = "  -7.34%  "
$ws.Range("D33").Value = "'1.093"
$ws.Range("D34").Value = "'1.721"
$ws.Range("E34").Value = "  -5.33%  "
$ws.Range("D35").Value = "'0.6438"
$ws.Range("E35").Value = "  -9.81%  "
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("D37").Value = "'2.731"
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("D38").Value = "'0.01743"
$ws.Range("E38").Value = "  -5.20%  "
$ws.Range("E39").Value = "  -3.82%  "
$ws.Range("D40").Value = "1.171.51"
$ws.Range("E40").Value = "  -5.00%  "
$ws.Range("D41").Value = "'0.8973"
$ws.Range("E41").Value = "  -6.04%  "
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").Value = "1.983.75"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").Value = "'100.75"
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("D45").Value = "'61.97"
$ws.Range("E45").Value = "  -5.13%  "
$ws.Range("D46").Value = "'0.00000000113"
$ws.Range("E46").Value = "  -3.68%  "
$ws.Range("D47").Value = "'0.05600"
$ws.Range("E47").Value = "  -0.95%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.466"
$ws.Range("E48").Value = "  -4.97%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.583"
$ws.Range("E49").Value = "  -6.85%  "
$ws.Range("D50").Value = "'0.4554"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").Value = "'0.3641"
$ws.Range("E51").Value = "  -6.27%  "
